$wb = $excel.ActiveWorkbook

# Sheet ALC (index 1), Row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 6541184.5
$ws.Range("I132").Value = 9263910
$ws.Range("J132").Value = 6644.933
$ws.Range("K132").Value = 27791730
$ws.Range("L132").Value = 19934.799
$ws.Range("M132").Value = -27789200
$ws.Range("N132").Value = -24994.799

# Sheet ALC (index 1), Row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 20833840
$ws.Range("I135").Value = 210.19048
$ws.Range("J135").Value = 166669250
$ws.Range("K135").Value = 1891.71432
$ws.Range("L135").Value = 1500023250
$ws.Range("M135").Value = 643.28568

# Sheet ALC (index 1), Row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1052.9565
$ws.Range("I137").Value = 714.2195
$ws.Range("J137").Value = 1548.9642
$ws.Range("K137").Value = 2142.6585
$ws.Range("L137").Value = 4646.892599999999
$ws.Range("M137").Value = 407.3415
$ws.Range("N137").Value = -9746.892599999999

# Sheet ALC (index 1), Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1351.7576
$ws.Range("I138").Value = 901.4054
$ws.Range("J138").Value = 1620.5161
$ws.Range("K138").Value = 2704.2162
$ws.Range("L138").Value = 4861.5483
$ws.Range("M138").Value = 2435.7838
$ws.Range("N138").Value = -15141.5483

# Sheet ALC (index 1), Row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 435.6389
$ws.Range("I141").Value = 435.6389
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1306.9167
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = $null
$ws.Range("N141").Value = 3873.0833

# Sheet ARM (index 2), Row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3973.175
$ws.Range("I32").Value = 3640.9412
$ws.Range("J32").Value = 5855.8335
$ws.Range("K32").Value = 3640.9412
$ws.Range("L32").Value = 5855.8335
$ws.Range("M32").Value = -3353.9412
$ws.Range("N32").Value = -6429.8335

# Sheet ARM (index 2), Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1279.1
$ws.Range("I45").Value = 1325.2858
$ws.Range("J45").Value = 1171.3334
$ws.Range("K45").Value = 1325.2858
$ws.Range("L45").Value = 1171.3334
$ws.Range("M45").Value = -948.2858000000001

# Sheet ARM (index 2), Row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1256.2903
$ws.Range("I74").Value = 931.6087
$ws.Range("J74").Value = 2189.75
$ws.Range("K74").Value = 931.6087
$ws.Range("L74").Value = 2189.75
$ws.Range("M74").Value = -57.6087
$ws.Range("N74").Value = -3937.75

# Sheet ARM (index 2), Row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1256.2903
$ws.Range("I77").Value = 931.6087
$ws.Range("J77").Value = 2189.75
$ws.Range("K77").Value = 4658.0435
$ws.Range("L77").Value = 10948.75
$ws.Range("M77").Value = -290.0434999999998
$ws.Range("N77").Value = -19684.75

# Sheet ARM (index 2), Row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1725.5
$ws.Range("I122").Value = 2075
$ws.Range("J122").Value = 910
$ws.Range("K122").Value = 6225
$ws.Range("L122").Value = 2730
$ws.Range("M122").Value = -3775
$ws.Range("N122").Value = -7630

# Sheet ARM (index 2), Row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1333.8909
$ws.Range("I132").Value = 1128.8379
$ws.Range("J132").Value = 1755.3889
$ws.Range("K132").Value = 3386.5137
$ws.Range("L132").Value = 5266.1667
$ws.Range("M132").Value = -856.5137
$ws.Range("N132").Value = -10326.1667

# Sheet BSM (index 3), Row 13
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = $null
$ws.Range("N13").Value = 0

# Sheet BSM (index 3), Row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 8333852
$ws.Range("I94").Value = 13158117
$ws.Range("J94").Value = 1029.6364
$ws.Range("K94").Value = 13158117
$ws.Range("L94").Value = 1029.6364
$ws.Range("M94").Value = -13157666
$ws.Range("N94").Value = -1931.6364

# Sheet CRP (index 4), Row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858910
$ws.Range("I16").Value = 166668220
$ws.Range("J16").Value = 3000
$ws.Range("K16").Value = 166668220
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -166667933

# Sheet CRP (index 4), Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1802.7333
$ws.Range("I31").Value = 1655.5526
$ws.Range("J31").Value = 2601.7144
$ws.Range("K31").Value = 1655.5526
$ws.Range("L31").Value = 2601.7144
$ws.Range("M31").Value = -1360.5526
$ws.Range("N31").Value = -3191.7144

# Sheet CRP (index 4), Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1802.7333
$ws.Range("I34").Value = 1655.5526
$ws.Range("J34").Value = 2601.7144
$ws.Range("K34").Value = 1655.5526
$ws.Range("L34").Value = 2601.7144
$ws.Range("M34").Value = -1453.5526
$ws.Range("N34").Value = -3005.7144

# Sheet CRP (index 4), Row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 911.129
$ws.Range("I58").Value = 804.3333
$ws.Range("J58").Value = 1632
$ws.Range("K58").Value = 804.3333
$ws.Range("L58").Value = 1632
$ws.Range("M58").Value = -601.3333
$ws.Range("N58").Value = -2038

# Sheet CRP (index 4), Row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 142858910
$ws.Range("I113").Value = 166668220
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 166668220
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -166666050

# Sheet CRP (index 4), Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6654.16
$ws.Range("I132").Value = 8113.4443
$ws.Range("J132").Value = 2901.7144
$ws.Range("K132").Value = 24340.3329
$ws.Range("L132").Value = 8705.143199999999
$ws.Range("M132").Value = -21810.3329
$ws.Range("N132").Value = -13765.1432

# Sheet CRP (index 4), Row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 20001828
$ws.Range("I134").Value = 1985.4286
$ws.Range("J134").Value = 125001000
$ws.Range("K134").Value = 5956.2858
$ws.Range("L134").Value = 375003000
$ws.Range("M134").Value = -3421.2858
$ws.Range("N134").Value = -375008070

# Sheet CRP (index 4), Row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 911.129
$ws.Range("I136").Value = 804.3333
$ws.Range("J136").Value = 1632
$ws.Range("K136").Value = 2412.9999
$ws.Range("L136").Value = 4896
$ws.Range("M136").Value = 137.0001000000002
$ws.Range("N136").Value = -9996

# Sheet CRP (index 4), Row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 29656.584
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 29656.584
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 29656.584
$ws.Range("N141").Value = -40016.584

# Sheet CUL (index 5), Row 55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2384.1667
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 3326.25
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 9978.75
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -10332.75

# Sheet CUL (index 5), Row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 619.4091
$ws.Range("I113").Value = 460.33334
$ws.Range("J113").Value = 660.3143
$ws.Range("K113").Value = 1381.00002
$ws.Range("L113").Value = 1980.9429
$ws.Range("M113").Value = 788.9999800000001
$ws.Range("N113").Value = -6320.9429

# Sheet CUL (index 5), Row 120
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 9128.375
$ws.Range("I120").Value = 2499.5
$ws.Range("J120").Value = 11338
$ws.Range("K120").Value = 7498.5
$ws.Range("L120").Value = 34014
$ws.Range("M120").Value = -2660.5
$ws.Range("N120").Value = -43690

# Sheet CUL (index 5), Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 25001350
$ws.Range("I131").Value = 100000560
$ws.Range("J131").Value = 1615
$ws.Range("K131").Value = 300001680
$ws.Range("L131").Value = 4845
$ws.Range("M131").Value = -299996640
$ws.Range("N131").Value = -14925

# Sheet CUL (index 5), Row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1430.7142
$ws.Range("I132").Value = 1053.5834
$ws.Range("J132").Value = 1933.5555
$ws.Range("K132").Value = 9482.250599999999
$ws.Range("L132").Value = 17401.9995
$ws.Range("M132").Value = -6952.250599999999
$ws.Range("N132").Value = -22461.9995

# Sheet CUL (index 5), Row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 20838760
$ws.Range("I137").Value = 46876784
$ws.Range("J137").Value = 8342.15
$ws.Range("K137").Value = 140630352
$ws.Range("L137").Value = 25026.45
$ws.Range("M137").Value = -140625252
$ws.Range("N137").Value = -35226.45

# Sheet CUL (index 5), Row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 24997.6
$ws.Range("I140").Value = 55188.316
$ws.Range("J140").Value = 2935.1538
$ws.Range("K140").Value = 165564.948
$ws.Range("L140").Value = 8805.4614
$ws.Range("M140").Value = -160384.948
$ws.Range("N140").Value = -19165.4614

# Sheet GSM (index 6), Row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 17858030
$ws.Range("I102").Value = 25000842
$ws.Range("J102").Value = 1003.5
$ws.Range("K102").Value = 25000842
$ws.Range("L102").Value = 1003.5
$ws.Range("M102").Value = -24999220
$ws.Range("N102").Value = -4247.5

# Sheet GSM (index 6), Row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2463.3044
$ws.Range("I132").Value = 1922.375
$ws.Range("J132").Value = 3699.7144
$ws.Range("K132").Value = 5767.125
$ws.Range("L132").Value = 11099.1432
$ws.Range("M132").Value = -3237.125

# Sheet LTW (index 7), Row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1990.7
$ws.Range("I7").Value = 1863.5
$ws.Range("J7").Value = 2499.5
$ws.Range("K7").Value = 1863.5
$ws.Range("L7").Value = 2499.5
$ws.Range("M7").Value = -1751.5

# Sheet LTW (index 7), Row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 1990.7
$ws.Range("I126").Value = 1863.5
$ws.Range("J126").Value = 2499.5
$ws.Range("K126").Value = 5590.5
$ws.Range("L126").Value = 7498.5
$ws.Range("M126").Value = -3120.5

# Sheet LTW (index 7), Row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 35840.332
$ws.Range("I132").Value = 2433.6667
$ws.Range("J132").Value = 50157.477
$ws.Range("K132").Value = 7301.000100000001
$ws.Range("L132").Value = 150472.431
$ws.Range("M132").Value = -4771.000100000001
$ws.Range("N132").Value = -155532.431

# Sheet LTW (index 7), Row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5315.8887
$ws.Range("I136").Value = 6780.4736
$ws.Range("J136").Value = 1837.5
$ws.Range("K136").Value = 20341.4208
$ws.Range("L136").Value = 5512.5
$ws.Range("M136").Value = -17791.4208
$ws.Range("N136").Value = -10612.5

# Sheet WVR (index 8), Row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 308.08694
$ws.Range("I113").Value = 190.63637
$ws.Range("J113").Value = 415.75
$ws.Range("K113").Value = 571.9091100000001
$ws.Range("L113").Value = 1247.25
$ws.Range("M113").Value = 1598.09089

# Sheet WVR (index 8), Row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1579.1578
$ws.Range("I132").Value = 1706
$ws.Range("J132").Value = 1486.909
$ws.Range("K132").Value = 5118
$ws.Range("L132").Value = 4460.727000000001
$ws.Range("M132").Value = -2588
$ws.Range("N132").Value = -9520.727000000001

# Sheet WVR (index 8), Row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 68857.5
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 68857.5
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 68857.5
$ws.Range("N135").Value = -78997.5

# Sheet WVR (index 8), Row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 759.8
$ws.Range("I136").Value = 801.3333
$ws.Range("J136").Value = 697.5
$ws.Range("K136").Value = 2403.9999
$ws.Range("L136").Value = 2092.5
$ws.Range("M136").Value = 146.0001000000002
$ws.Range("N136").Value = -7192.5
